# Auto-generated edit script: applies scheduled-runner value updates
# to the Seraph_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Cells.Item(4, 8).Value = 250.375  # H4: 281.57144 -> 250.375
$ws.Cells.Item(4, 9).Value = 143.14285  # I4: 161.66667 -> 143.14285
$ws.Cells.Item(4, 11).Value = 143.14285  # K4: 161.66667 -> 143.14285
$ws.Cells.Item(4, 13).Value = -29.14285000000001  # M4: -47.66667000000001 -> -29.14285000000001

# Row 11
$ws.Cells.Item(11, 8).Value = 257.7  # H11: 254.18182 -> 257.7
$ws.Cells.Item(11, 9).Value = 257.7  # I11: 254.18182 -> 257.7
$ws.Cells.Item(11, 11).Value = 257.7  # K11: 254.18182 -> 257.7
$ws.Cells.Item(11, 13).Value = -117.7  # M11: -114.18182 -> -117.7

# Row 64
$ws.Cells.Item(64, 8).Value = 4299.6665  # H64: 4371.2856 -> 4299.6665
$ws.Cells.Item(64, 9).Value = 4519.6  # I64: 6000 -> 4519.6
$ws.Cells.Item(64, 10).Value = 3200  # J64: 3149.75 -> 3200
$ws.Cells.Item(64, 11).Value = 4519.6  # K64: 6000 -> 4519.6
$ws.Cells.Item(64, 12).Value = 3200  # L64: 3149.75 -> 3200
$ws.Cells.Item(64, 13).Value = -4271.6  # M64: -5752 -> -4271.6
$ws.Cells.Item(64, 14).Value = -3696  # N64: -3645.75 -> -3696

# Row 67
$ws.Cells.Item(67, 8).Value = 4299.6665  # H67: 4371.2856 -> 4299.6665
$ws.Cells.Item(67, 9).Value = 4519.6  # I67: 6000 -> 4519.6
$ws.Cells.Item(67, 10).Value = 3200  # J67: 3149.75 -> 3200
$ws.Cells.Item(67, 11).Value = 4519.6  # K67: 6000 -> 4519.6
$ws.Cells.Item(67, 12).Value = 3200  # L67: 3149.75 -> 3200
$ws.Cells.Item(67, 13).Value = -3661.6  # M67: -5142 -> -3661.6
$ws.Cells.Item(67, 14).Value = -4916  # N67: -4865.75 -> -4916

# Row 137
$ws.Cells.Item(137, 8).Value = 2560.4614  # H137: 2527.6428 -> 2560.4614
$ws.Cells.Item(137, 9).Value = 1115.6666  # I137: 1133 -> 1115.6666
$ws.Cells.Item(137, 10).Value = 3798.8572  # J137: 3573.625 -> 3798.8572
$ws.Cells.Item(137, 11).Value = 3346.9998  # K137: 3399 -> 3346.9998
$ws.Cells.Item(137, 12).Value = 11396.5716  # L137: 10720.875 -> 11396.5716
$ws.Cells.Item(137, 13).Value = -796.9998000000001  # M137: -849 -> -796.9998000000001
$ws.Cells.Item(137, 14).Value = -16496.5716  # N137: -15820.875 -> -16496.5716

# Row 138
$ws.Cells.Item(138, 8).Value = 4228.0303  # H138: 4384.0864 -> 4228.0303
$ws.Cells.Item(138, 10).Value = 4544.5  # J138: 4776.16 -> 4544.5
$ws.Cells.Item(138, 12).Value = 13633.5  # L138: 14328.48 -> 13633.5
$ws.Cells.Item(138, 14).Value = -23913.5  # N138: -24608.48 -> -23913.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 1499.5  # H61: 0 -> 1499.5
$ws.Cells.Item(61, 9).Value = 1500  # I61: 0 -> 1500
$ws.Cells.Item(61, 10).Value = 1499  # J61: 0 -> 1499
$ws.Cells.Item(61, 11).Value = 1500  # K61: 0 -> 1500
$ws.Cells.Item(61, 12).Value = 1499  # L61: 0 -> 1499
$ws.Cells.Item(61, 13).Value = -1288  # M61: None -> -1288
$ws.Cells.Item(61, 14).Value = -1923  # N61: None -> -1923

# Row 82
$ws.Cells.Item(82, 8).Value = 39044.668  # H82: 0 -> 39044.668
$ws.Cells.Item(82, 10).Value = 39044.668  # J82: 0 -> 39044.668
$ws.Cells.Item(82, 12).Value = 39044.668  # L82: 0 -> 39044.668
$ws.Cells.Item(82, 14).Value = -39766.668  # N82: None -> -39766.668

# Row 85
$ws.Cells.Item(85, 8).Value = 39044.668  # H85: 0 -> 39044.668
$ws.Cells.Item(85, 10).Value = 39044.668  # J85: 0 -> 39044.668
$ws.Cells.Item(85, 12).Value = 39044.668  # L85: 0 -> 39044.668
$ws.Cells.Item(85, 14).Value = -41540.668  # N85: None -> -41540.668

# Row 122
$ws.Cells.Item(122, 8).Value = 670403.9  # H122: 718082.2 -> 670403.9
$ws.Cells.Item(122, 9).Value = 1003105.8  # I122: 1114239 -> 1003105.8
$ws.Cells.Item(122, 11).Value = 3009317.4  # K122: 3342717 -> 3009317.4
$ws.Cells.Item(122, 13).Value = -3006867.4  # M122: -3340267 -> -3006867.4

# Row 132
$ws.Cells.Item(132, 8).Value = 1456.1621  # H132: 1435.2894 -> 1456.1621
$ws.Cells.Item(132, 9).Value = 1482.3143  # I132: 1459.5555 -> 1482.3143
$ws.Cells.Item(132, 11).Value = 4446.9429  # K132: 4378.666499999999 -> 4446.9429
$ws.Cells.Item(132, 13).Value = -1916.9429  # M132: -1848.666499999999 -> -1916.9429

# Row 136
$ws.Cells.Item(136, 8).Value = 1499.5  # H136: 0 -> 1499.5
$ws.Cells.Item(136, 9).Value = 1500  # I136: 0 -> 1500
$ws.Cells.Item(136, 10).Value = 1499  # J136: 0 -> 1499
$ws.Cells.Item(136, 11).Value = 4500  # K136: 0 -> 4500
$ws.Cells.Item(136, 12).Value = 4497  # L136: 0 -> 4497
$ws.Cells.Item(136, 13).Value = -1950  # M136: None -> -1950
$ws.Cells.Item(136, 14).Value = -9597  # N136: None -> -9597

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1355.7858  # H20: 1421.6154 -> 1355.7858
$ws.Cells.Item(20, 9).Value = 1210.2222  # I20: 1299 -> 1210.2222
$ws.Cells.Item(20, 11).Value = 1210.2222  # K20: 1299 -> 1210.2222
$ws.Cells.Item(20, 13).Value = -963.2221999999999  # M20: -1052 -> -963.2221999999999

# Row 75
$ws.Cells.Item(75, 8).Value = 50398.8  # H75: 43703 -> 50398.8
$ws.Cells.Item(75, 9).Value = 12000  # I75: 11111 -> 12000
$ws.Cells.Item(75, 10).Value = 59998.5  # J75: 59999 -> 59998.5
$ws.Cells.Item(75, 11).Value = 12000  # K75: 11111 -> 12000
$ws.Cells.Item(75, 12).Value = 59998.5  # L75: 59999 -> 59998.5
$ws.Cells.Item(75, 13).Value = -11064  # M75: -10175 -> -11064
$ws.Cells.Item(75, 14).Value = -61870.5  # N75: -61871 -> -61870.5

# Row 78
$ws.Cells.Item(78, 8).Value = 50398.8  # H78: 43703 -> 50398.8
$ws.Cells.Item(78, 9).Value = 12000  # I78: 11111 -> 12000
$ws.Cells.Item(78, 10).Value = 59998.5  # J78: 59999 -> 59998.5
$ws.Cells.Item(78, 11).Value = 36000  # K78: 33333 -> 36000
$ws.Cells.Item(78, 12).Value = 179995.5  # L78: 179997 -> 179995.5
$ws.Cells.Item(78, 13).Value = -31320  # M78: -28653 -> -31320
$ws.Cells.Item(78, 14).Value = -189355.5  # N78: -189357 -> -189355.5

# Row 122
$ws.Cells.Item(122, 8).Value = 424998  # H122: 424999.2 -> 424998
$ws.Cells.Item(122, 9).Value = 299997  # I122: 299999 -> 299997
$ws.Cells.Item(122, 10).Value = 508332  # J122: 508332.66 -> 508332
$ws.Cells.Item(122, 11).Value = 299997  # K122: 299999 -> 299997
$ws.Cells.Item(122, 12).Value = 508332  # L122: 508332.66 -> 508332
$ws.Cells.Item(122, 13).Value = -295097  # M122: -295099 -> -295097
$ws.Cells.Item(122, 14).Value = -518132  # N122: -518132.66 -> -518132

# Row 134
$ws.Cells.Item(134, 8).Value = 3429.08  # H134: 3472.0417 -> 3429.08
$ws.Cells.Item(134, 9).Value = 3447  # I134: 3492.6086 -> 3447
$ws.Cells.Item(134, 11).Value = 10341  # K134: 10477.8258 -> 10341
$ws.Cells.Item(134, 13).Value = -7806  # M134: -7942.825800000001 -> -7806

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Cells.Item(6, 8).Value = 2144907.2  # H6: 1667566.6 -> 2144907.2
$ws.Cells.Item(6, 9).Value = 3753512.5  # I6: 5000 -> 3753512.5
$ws.Cells.Item(6, 10).Value = 100  # J6: 2000080 -> 100
$ws.Cells.Item(6, 11).Value = 3753512.5  # K6: 5000 -> 3753512.5
$ws.Cells.Item(6, 12).Value = 100  # L6: 2000080 -> 100
$ws.Cells.Item(6, 13).Value = -3753399.5  # M6: -4887 -> -3753399.5
$ws.Cells.Item(6, 14).Value = -326  # N6: -2000306 -> -326

# Row 7
$ws.Cells.Item(7, 8).Value = 181.5  # H7: 138.91667 -> 181.5
$ws.Cells.Item(7, 9).Value = 235.33333  # I7: 227.85715 -> 235.33333
$ws.Cells.Item(7, 10).Value = 20  # J7: 14.4 -> 20
$ws.Cells.Item(7, 11).Value = 235.33333  # K7: 227.85715 -> 235.33333
$ws.Cells.Item(7, 12).Value = 20  # L7: 14.4 -> 20
$ws.Cells.Item(7, 13).Value = -122.33333  # M7: -114.85715 -> -122.33333
$ws.Cells.Item(7, 14).Value = -246  # N7: -240.4 -> -246

# Row 11
$ws.Cells.Item(11, 8).Value = 4666.6665  # H11: 1500 -> 4666.6665
$ws.Cells.Item(11, 9).Value = 3000  # I11: 0 -> 3000
$ws.Cells.Item(11, 10).Value = 8000  # J11: 1500 -> 8000
$ws.Cells.Item(11, 11).Value = 3000  # K11: 0 -> 3000
$ws.Cells.Item(11, 12).Value = 8000  # L11: 1500 -> 8000
$ws.Cells.Item(11, 13).Value = -2860  # M11: None -> -2860
$ws.Cells.Item(11, 14).Value = -8280  # N11: -1780 -> -8280

# Row 12
$ws.Cells.Item(12, 8).Value = 5500  # H12: 7874.75 -> 5500
$ws.Cells.Item(12, 9).Value = 10000  # I12: 7166.3335 -> 10000
$ws.Cells.Item(12, 10).Value = 1000  # J12: 10000 -> 1000
$ws.Cells.Item(12, 11).Value = 10000  # K12: 7166.3335 -> 10000
$ws.Cells.Item(12, 12).Value = 1000  # L12: 10000 -> 1000
$ws.Cells.Item(12, 13).Value = -9830  # M12: -6996.3335 -> -9830
$ws.Cells.Item(12, 14).Value = -1340  # N12: -10340 -> -1340

# Row 13
$ws.Cells.Item(13, 8).Value = 19999  # H13: 0 -> 19999
$ws.Cells.Item(13, 10).Value = 19999  # J13: 0 -> 19999
$ws.Cells.Item(13, 12).Value = 19999  # L13: 0 -> 19999
$ws.Cells.Item(13, 14).Value = -20277  # N13: None -> -20277

# Row 31
$ws.Cells.Item(31, 8).Value = 4486.8096  # H31: 4491.864 -> 4486.8096
$ws.Cells.Item(31, 9).Value = 2798.2  # I31: 2858.2 -> 2798.2
$ws.Cells.Item(31, 10).Value = 6021.909  # J31: 5853.25 -> 6021.909
$ws.Cells.Item(31, 11).Value = 2798.2  # K31: 2858.2 -> 2798.2
$ws.Cells.Item(31, 12).Value = 6021.909  # L31: 5853.25 -> 6021.909
$ws.Cells.Item(31, 13).Value = -2503.2  # M31: -2563.2 -> -2503.2
$ws.Cells.Item(31, 14).Value = -6611.909  # N31: -6443.25 -> -6611.909

# Row 34
$ws.Cells.Item(34, 8).Value = 4486.8096  # H34: 4491.864 -> 4486.8096
$ws.Cells.Item(34, 9).Value = 2798.2  # I34: 2858.2 -> 2798.2
$ws.Cells.Item(34, 10).Value = 6021.909  # J34: 5853.25 -> 6021.909
$ws.Cells.Item(34, 11).Value = 2798.2  # K34: 2858.2 -> 2798.2
$ws.Cells.Item(34, 12).Value = 6021.909  # L34: 5853.25 -> 6021.909
$ws.Cells.Item(34, 13).Value = -2596.2  # M34: -2656.2 -> -2596.2
$ws.Cells.Item(34, 14).Value = -6425.909  # N34: -6257.25 -> -6425.909

$ws = $wb.Worksheets.Item("CUL")
# Row 60
$ws.Cells.Item(60, 8).Value = 499.16666  # H60: 733 -> 499.16666
$ws.Cells.Item(60, 9).Value = 267.5  # I60: 618.25 -> 267.5
$ws.Cells.Item(60, 11).Value = 802.5  # K60: 1854.75 -> 802.5
$ws.Cells.Item(60, 13).Value = -551.5  # M60: -1603.75 -> -551.5

# Row 75
$ws.Cells.Item(75, 8).Value = 367.5  # H75: 372.14285 -> 367.5
$ws.Cells.Item(75, 9).Value = 391  # I75: 392.5 -> 391
$ws.Cells.Item(75, 11).Value = 1173  # K75: 1177.5 -> 1173
$ws.Cells.Item(75, 13).Value = -175  # M75: -179.5 -> -175

# Row 78
$ws.Cells.Item(78, 8).Value = 367.5  # H78: 372.14285 -> 367.5
$ws.Cells.Item(78, 9).Value = 391  # I78: 392.5 -> 391
$ws.Cells.Item(78, 11).Value = 3519  # K78: 3532.5 -> 3519
$ws.Cells.Item(78, 13).Value = 1473  # M78: 1459.5 -> 1473

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Cells.Item(15, 8).Value = 14950  # H15: 13300 -> 14950
$ws.Cells.Item(15, 10).Value = 14950  # J15: 13300 -> 14950
$ws.Cells.Item(15, 12).Value = 14950  # L15: 13300 -> 14950
$ws.Cells.Item(15, 14).Value = -15526  # N15: -13876 -> -15526

# Row 70
$ws.Cells.Item(70, 8).Value = 6855.7144  # H70: 7998 -> 6855.7144
$ws.Cells.Item(70, 9).Value = 4000  # I70: 0 -> 4000
$ws.Cells.Item(70, 11).Value = 4000  # K70: 0 -> 4000
$ws.Cells.Item(70, 13).Value = -3730  # M70: None -> -3730

# Row 73
$ws.Cells.Item(73, 8).Value = 6855.7144  # H73: 7998 -> 6855.7144
$ws.Cells.Item(73, 9).Value = 4000  # I73: 0 -> 4000
$ws.Cells.Item(73, 11).Value = 4000  # K73: 0 -> 4000
$ws.Cells.Item(73, 13).Value = -3064  # M73: None -> -3064

# Row 80
$ws.Cells.Item(80, 8).Value = 4373.6665  # H80: 4589.909 -> 4373.6665
$ws.Cells.Item(80, 9).Value = 3286.875  # I80: 3471.4285 -> 3286.875
$ws.Cells.Item(80, 11).Value = 3286.875  # K80: 3471.4285 -> 3286.875
$ws.Cells.Item(80, 13).Value = -2288.875  # M80: -2473.4285 -> -2288.875

# Row 81
$ws.Cells.Item(81, 8).Value = 14950  # H81: 13300 -> 14950
$ws.Cells.Item(81, 10).Value = 14950  # J81: 13300 -> 14950
$ws.Cells.Item(81, 12).Value = 14950  # L81: 13300 -> 14950
$ws.Cells.Item(81, 14).Value = -16946  # N81: -15296 -> -16946

# Row 83
$ws.Cells.Item(83, 8).Value = 4373.6665  # H83: 4589.909 -> 4373.6665
$ws.Cells.Item(83, 9).Value = 3286.875  # I83: 3471.4285 -> 3286.875
$ws.Cells.Item(83, 11).Value = 16434.375  # K83: 17357.1425 -> 16434.375
$ws.Cells.Item(83, 13).Value = -11442.375  # M83: -12365.1425 -> -11442.375

# Row 84
$ws.Cells.Item(84, 8).Value = 14950  # H84: 13300 -> 14950
$ws.Cells.Item(84, 10).Value = 14950  # J84: 13300 -> 14950
$ws.Cells.Item(84, 12).Value = 44850  # L84: 39900 -> 44850
$ws.Cells.Item(84, 14).Value = -54834  # N84: -49884 -> -54834

# Row 113
$ws.Cells.Item(113, 8).Value = 2565.2778  # H113: 2668.85 -> 2565.2778
$ws.Cells.Item(113, 9).Value = 1347.9166  # I113: 1424 -> 1347.9166
$ws.Cells.Item(113, 10).Value = 5000  # J113: 4190.3335 -> 5000
$ws.Cells.Item(113, 11).Value = 1347.9166  # K113: 1424 -> 1347.9166
$ws.Cells.Item(113, 12).Value = 5000  # L113: 4190.3335 -> 5000
$ws.Cells.Item(113, 13).Value = 822.0834  # M113: 746 -> 822.0834
$ws.Cells.Item(113, 14).Value = -9340  # N113: -8530.333500000001 -> -9340

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Cells.Item(55, 8).Value = 352.95  # H55: 340.85715 -> 352.95
$ws.Cells.Item(55, 10).Value = 416.4  # J55: 363.5 -> 416.4
$ws.Cells.Item(55, 12).Value = 416.4  # L55: 363.5 -> 416.4
$ws.Cells.Item(55, 14).Value = -762.4  # N55: -709.5 -> -762.4

# Row 61
$ws.Cells.Item(61, 8).Value = 2631.8235  # H61: 2574.5 -> 2631.8235
$ws.Cells.Item(61, 9).Value = 2133.7693  # I61: 2168.4614 -> 2133.7693
$ws.Cells.Item(61, 10).Value = 4250.5  # J61: 4334 -> 4250.5
$ws.Cells.Item(61, 11).Value = 2133.7693  # K61: 2168.4614 -> 2133.7693
$ws.Cells.Item(61, 12).Value = 4250.5  # L61: 4334 -> 4250.5
$ws.Cells.Item(61, 13).Value = -1931.7693  # M61: -1966.4614 -> -1931.7693
$ws.Cells.Item(61, 14).Value = -4654.5  # N61: -4738 -> -4654.5

# Row 82
$ws.Cells.Item(82, 8).Value = 2356.2856  # H82: 2762.4546 -> 2356.2856
$ws.Cells.Item(82, 9).Value = 2320.889  # I82: 3417.4 -> 2320.889
$ws.Cells.Item(82, 10).Value = 2420  # J82: 2216.6667 -> 2420
$ws.Cells.Item(82, 11).Value = 2320.889  # K82: 3417.4 -> 2320.889
$ws.Cells.Item(82, 12).Value = 2420  # L82: 2216.6667 -> 2420
$ws.Cells.Item(82, 13).Value = -1959.889  # M82: -3056.4 -> -1959.889
$ws.Cells.Item(82, 14).Value = -3142  # N82: -2938.6667 -> -3142

# Row 85
$ws.Cells.Item(85, 8).Value = 2356.2856  # H85: 2762.4546 -> 2356.2856
$ws.Cells.Item(85, 9).Value = 2320.889  # I85: 3417.4 -> 2320.889
$ws.Cells.Item(85, 10).Value = 2420  # J85: 2216.6667 -> 2420
$ws.Cells.Item(85, 11).Value = 2320.889  # K85: 3417.4 -> 2320.889
$ws.Cells.Item(85, 12).Value = 2420  # L85: 2216.6667 -> 2420
$ws.Cells.Item(85, 13).Value = -1072.889  # M85: -2169.4 -> -1072.889
$ws.Cells.Item(85, 14).Value = -4916  # N85: -4712.6667 -> -4916

# Row 113
$ws.Cells.Item(113, 8).Value = 2631.8235  # H113: 2574.5 -> 2631.8235
$ws.Cells.Item(113, 9).Value = 2133.7693  # I113: 2168.4614 -> 2133.7693
$ws.Cells.Item(113, 10).Value = 4250.5  # J113: 4334 -> 4250.5
$ws.Cells.Item(113, 11).Value = 2133.7693  # K113: 2168.4614 -> 2133.7693
$ws.Cells.Item(113, 12).Value = 4250.5  # L113: 4334 -> 4250.5
$ws.Cells.Item(113, 13).Value = 36.23070000000007  # M113: 1.53859999999986 -> 36.23070000000007
$ws.Cells.Item(113, 14).Value = -8590.5  # N113: -8674 -> -8590.5

# Row 132
$ws.Cells.Item(132, 8).Value = 5244.6665  # H132: 5263.1143 -> 5244.6665
$ws.Cells.Item(132, 9).Value = 4885.8696  # I132: 4898.909 -> 4885.8696
$ws.Cells.Item(132, 11).Value = 14657.6088  # K132: 14696.727 -> 14657.6088
$ws.Cells.Item(132, 13).Value = -12127.6088  # M132: -12166.727 -> -12127.6088

$ws = $wb.Worksheets.Item("WVR")
# Row 11
$ws.Cells.Item(11, 8).Value = 0  # H11: 25002.5 -> 0
$ws.Cells.Item(11, 10).Value = 0  # J11: 25002.5 -> 0
$ws.Cells.Item(11, 12).Value = 0  # L11: 25002.5 -> 0
$ws.Cells.Item(11, 14).ClearContents()  # N11: -25286.5 -> (removed)

# Row 132
$ws.Cells.Item(132, 8).Value = 1601.3636  # H132: 1611.7 -> 1601.3636
$ws.Cells.Item(132, 10).Value = 1582.3334  # J132: 1624.5 -> 1582.3334
$ws.Cells.Item(132, 12).Value = 4747.0002  # L132: 4873.5 -> 4747.0002
$ws.Cells.Item(132, 14).Value = -9807.0002  # N132: -9933.5 -> -9807.0002

Write-Output "Applied 43 cell-range updates across 8 sheets."